$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{A="WC48 P5F"; B="Etiquetadora"; C="2024-06-10"; D="09:53:28"; E="Mañana"; F="09:53:29"; G="0:00:01"; H="N/A"}
    @{A="WC48 P5F"; B="Cámara no detecta foams"; C="2024-06-10"; D="09:53:30"; E="Mañana"; F="09:53:31"; G="0:00:01"; H="N/A"}
    @{A="WC48 P5F"; B="Cámara no detecta foam derecho"; C="2024-06-10"; D="09:53:38"; E="Mañana"; F="09:53:38"; G="0:00:00"; H="N/A"}
    @{A="WC47 NACP"; B="Fallo tolva"; C="2024-06-10"; D="09:57:07"; E="Mañana"; F="09:57:08"; G="0:00:01"; H="-0.00 minutos"}
    @{A="WC47 NACP"; B="Ascensor no sube"; C="2024-06-10"; D="09:57:10"; E="Mañana"; F="09:57:12"; G="0:00:02"; H="0.01 minutos"}
    @{A="WC47 NACP"; B="Etiquetadora"; C="2024-06-10"; D="09:57:13"; E="Mañana"; F="09:57:15"; G="0:00:02"; H="0.02 minutos"}
    @{A="WC47 NACP"; B="Fallo en paletizador"; C="2024-06-10"; D="09:57:30"; E="Mañana"; F="09:57:31"; G="0:00:01"; H="0.05 minutos"}
    @{A="WC47 NACP"; B="No coge placa"; C="2024-06-10"; D="09:58:21"; E="Mañana"; F="09:59:09"; G="0:00:48"; H="0.12 minutos"}
    @{A="WC47 NACP"; B="Fallo tornillo"; C="2024-06-10"; D="09:59:30"; E="Mañana"; F="09:59:36"; G="0:00:06"; H="0.20 minutos"}
    @{A="WC47 NACP"; B="Fallo tolva"; C="2024-06-10"; D="09:59:49"; E="Mañana"; F="09:59:51"; G="0:00:02"; H="0.19 minutos"}
)

# Column C holds dates formatted as plain text (e.g. "2024-06-10"); force
# text format up front so Excel doesn't auto-convert them to date serials.
$ws.Range("C149:C158").NumberFormat = "@"

$r = 149
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r++
}
